$d = $word.ActiveDocument

$map = @(
    @("94×91=8554", "84×97=8148"),
    @("67×42=2814", "98×57=5586"),
    @("82×67=5494", "24×27=648"),
    @("89×31=2759", "91×45=4095"),
    @("50×15=750",  "90×87=7830"),
    @("64×55=3520", "16×40=640"),
    @("91×21=1911", "98×92=9016"),
    @("41×39=1599", "59×57=3363"),
    @("99×51=5049", "30×65=1950"),
    @("91×56=5096", "80×47=3760"),
    @("20×63=1260", "18×52=936"),
    @("27×63=1701", "37×76=2812"),
    @("62×82=5084", "32×91=2912"),
    @("96×94=9024", "74×42=3108"),
    @("46×16=736",  "82×84=6888"),
    @("72×29=2088", "39×19=741"),
    @("46×96=4416", "84×98=8232"),
    @("79×61=4819", "55×52=2860"),
    @("92×69=6348", "59×13=767"),
    @("51×40=2040", "11×65=715"),
    @("49×71=3479", "16×50=800"),
    @("21×40=840",  "63×61=3843"),
    @("77×43=3311", "70×72=5040"),
    @("18×20=360",  "76×95=7220"),
    @("47×94=4418", "87×34=2958")
)

foreach ($pair in $map) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
